$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# 1. Rename the header text of column E from "MW.Region.code" to "AMW.Region.code"
$ws.Range("E1").Value2 = "AMW.Region.code"

# 2. Fill column F for rows 206-240 with the same value as column E ("-"),
#    matching the pattern already used elsewhere in the sheet (e.g. row 18).
for ($r = 206; $r -le 240; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $fCell.Value2 = $eCell.Value2
}

# 3. Update the active selection to reflect F206:F240 with active cell F206
$ws.Range("F206:F240").Select()
